$d = $word.ActiveDocument

$replacements = @(
    @("2025-01-29 Wednesday", "2025-01-30 Thursday"),
    @("53÷7=7, 4", "79÷3=26, 1"),
    @("25÷8=3, 1", "25÷9=2, 7"),
    @("29÷6=4, 5", "11÷7=1, 4"),
    @("28÷7=4, 0", "95÷9=10, 5"),
    @("58÷7=8, 2", "28÷9=3, 1"),
    @("45÷6=7, 3", "63÷2=31, 1"),
    @("98÷2=49, 0", "33÷8=4, 1"),
    @("93÷2=46, 1", "30÷6=5, 0"),
    @("20÷6=3, 2", "60÷7=8, 4"),
    @("47÷7=6, 5", "63÷8=7, 7"),
    @("22÷4=5, 2", "16÷6=2, 4"),
    @("70÷2=35, 0", "21÷4=5, 1"),
    @("94÷9=10, 4", "62÷7=8, 6"),
    @("90÷6=15, 0", "24÷8=3, 0"),
    @("59÷4=14, 3", "33÷7=4, 5"),
    @("65÷5=13, 0", "50÷6=8, 2"),
    @("56÷7=8, 0", "46÷2=23, 0"),
    @("61÷5=12, 1", "63÷8=7, 7"),
    @("82÷7=11, 5", "18÷8=2, 2"),
    @("99÷5=19, 4", "76÷6=12, 4"),
    @("36÷4=9, 0", "77÷6=12, 5"),
    @("48÷7=6, 6", "17÷9=1, 8"),
    @("19÷9=2, 1", "24÷7=3, 3"),
    @("40÷3=13, 1", "59÷3=19, 2"),
    @("80÷5=16, 0", "70÷5=14, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
